$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 20
# from 2023-10-09 (45208) to 2023-10-13 (45212)
$ws.Range("C2:C20").Value = 45212
